$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row above row 2 (shifts existing rows 2-35 down to 3-36,
# and auto-adjusts chart/formula references accordingly).
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row 2 with the new data series entry.
$ws.Range("A2").Value = "Python3 (second year)"
$ws.Range("B2").Value = 11
$ws.Range("C2").Value = 13

$ws.Range("M26").Select()
